# working_hours.xlsx — add a new time-tracking entry for 2014-07-28.
#
# The sheet lists one work interval per row (year/month/day/start/end) with
# helper columns F (minutes) / G (hours), followed by a blank spacer row and
# three summary rows (sum [min], sum [h], sum [working weeks]). We insert a
# new data row right above the old spacer/summary block and backfill the
# summary formulas so the totals include it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing spacer row (old row 161) and the three summary rows
# (old rows 162-164) down by one to make room for the new data row.
$ws.Rows("161:161").Insert()

# New entry: 2014-07-28, 08:30 -> 12:00 (3.5h / 210min), same pattern/styles
# as every other data row.
$ws.Range("A161").Value = 2014
$ws.Range("B161").Value = 7
$ws.Range("C161").Value = 28
$ws.Range("D161").Value = 0.35416666666666669
$ws.Range("E161").Value = 0.5
$ws.Range("F161").Formula = "=(E161-D161)*24*60"
$ws.Range("G161").Formula = "=F161/60"

# The summary formulas (now on rows 163-165) still pointed at the old data
# range; extend the sum to include the newly inserted row.
$ws.Range("F163").Formula = "=SUM(F2:F161)"

# Reflect where the user ended up after typing the new row.
$ws.Range("A162").Select() | Out-Null
